# #12 Web page added on slides
#
# Inserts a new "Web page" slide right before the existing "Strategy
# learning" slide (previously slide 16), using the same "Title and
# Content" layout as its neighbours. This pushes "Strategy learning",
# "Result", "Conclusion", "Future work", etc. down by one position.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(16, 2)

$title = $newSlide.Shapes.Item(1)
$title.Name = "Titre 1"
$title.TextFrame.TextRange.Text = "Web page"

$body = $newSlide.Shapes.Item(2)
$body.Name = "Espace réservé du contenu 2"
$body.TextFrame.TextRange.Text = "The web page will show us all the strategies SAI record during the day`rExplain with a schema"
